$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.296.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9983"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.94%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9959"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4369"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +14.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3536"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.50"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.159"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07464"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.99"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9973"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.288"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.306"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.816.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06679"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9978"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.472"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.299.34"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.10"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.375"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.483"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.15"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.026.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.308"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -11.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.90"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.060"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.982"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09363"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6790"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02371"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2173"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.492"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.220"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.283"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9959"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6170"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.871"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.049"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.175"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07110"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.31%  "

# Row 38: swap to Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06277"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.14%  "

# Row 39: swap to InternetComputer(DFINITY)
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.231"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.48%  "
